# Apply "feat: add 2022-Q3 data":
#   1. Insert a new worksheet "2022-Q3" right after "总计", seeded from the
#      "2022-Q2" sheet's layout/styles, then overwrite its fund figures.
#   2. Update the "总计" summary sheet: push the existing quarters down one
#      row and insert the new 2022-Q3 totals at the top of the table.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    # Force a numeric-looking string (e.g. "5.89") to be stored as text
    # instead of being auto-coerced to a number, then drop back to the
    # default "Normal" style so no stray number-format style lingers on
    # the cell.
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q3" sheet right after "总计" by copying
# the "2022-Q2" sheet (keeps header row / column styles identical).
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($null, $wsTotal)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# Row 2: 009562 (unchanged fund code)
Set-TextValue $wsQ3.Range("C2") "工银全球股票（QDII）美元"
Set-TextValue $wsQ3.Range("D2") "5.89"
Set-TextValue $wsQ3.Range("E2") "93.72"
Set-TextValue $wsQ3.Range("F2") "2.27"
Set-TextValue $wsQ3.Range("G2") "0.1337"
$wsQ3.Range("H2").Value2 = 6

# Row 3: 009563 (unchanged fund code)
Set-TextValue $wsQ3.Range("C3") "工银全球股票（QDII）港币"
Set-TextValue $wsQ3.Range("D3") "5.89"
Set-TextValue $wsQ3.Range("E3") "93.72"
Set-TextValue $wsQ3.Range("F3") "2.27"
Set-TextValue $wsQ3.Range("G3") "0.1337"
$wsQ3.Range("H3").Value2 = 6

# Row 4: 486001 (unchanged fund code)
Set-TextValue $wsQ3.Range("C4") "工银瑞信中国机会全球配置股票（QDII）人民币"
Set-TextValue $wsQ3.Range("D4") "5.89"
Set-TextValue $wsQ3.Range("E4") "93.72"
Set-TextValue $wsQ3.Range("F4") "2.27"
Set-TextValue $wsQ3.Range("G4") "0.1337"
$wsQ3.Range("H4").Value2 = 6

# Row 5: 539002 (unchanged fund code)
Set-TextValue $wsQ3.Range("C5") "建信新兴市场优选混合（QDII）"
Set-TextValue $wsQ3.Range("D5") "0.14"
Set-TextValue $wsQ3.Range("E5") "81.57"
Set-TextValue $wsQ3.Range("F5") "2.31"
Set-TextValue $wsQ3.Range("G5") "0.0032"
$wsQ3.Range("H5").Value2 = 10

# ---------------------------------------------------------------------
# Step 2: update the "总计" sheet - shift rows 2..8 down to 3..9, then
# write the new 2022-Q3 row at row 2 and fix up the newly exposed A9.
# ---------------------------------------------------------------------
for ($r = 8; $r -ge 2; $r--) {
    $nr = $r + 1
    $wsTotal.Range("B$nr").Value2 = $wsTotal.Range("B$r").Value2
    $wsTotal.Range("C$nr").Value2 = $wsTotal.Range("C$r").Value2
    $wsTotal.Range("D$nr").Value2 = $wsTotal.Range("D$r").Value2
}

$wsTotal.Range("A8").Copy($wsTotal.Range("A9"))
$wsTotal.Range("A9").Value2 = 7

$wsTotal.Range("B2").Value2 = "2022-Q3"
$wsTotal.Range("C2").Value2 = 4
$wsTotal.Range("D2").Value2 = 0.4
